# Update the multiplication problems in the table to new values.
$d = $word.ActiveDocument

$replacements = @(
    @{old="470×9="; new="395×8="},
    @{old="676×8="; new="699×9="},
    @{old="786×5="; new="208×9="},
    @{old="290×6="; new="125×5="},
    @{old="649×8="; new="125×9="},
    @{old="156×5="; new="689×9="},
    @{old="491×6="; new="193×2="},
    @{old="795×5="; new="655×6="},
    @{old="992×5="; new="789×2="},
    @{old="290×8="; new="648×6="},
    @{old="498×5="; new="826×8="},
    @{old="670×3="; new="187×8="},
    @{old="730×7="; new="771×8="},
    @{old="467×2="; new="853×4="},
    @{old="206×8="; new="654×6="},
    @{old="546×5="; new="829×6="},
    @{old="748×9="; new="382×4="},
    @{old="972×6="; new="989×3="},
    @{old="232×5="; new="617×7="},
    @{old="869×3="; new="242×2="},
    @{old="120×3="; new="675×9="},
    @{old="300×6="; new="452×7="},
    @{old="633×7="; new="135×7="},
    @{old="920×6="; new="969×8="},
    @{old="960×7="; new="243×4="}
)

foreach ($r in $replacements) {
    $found = $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
    Write-Host "Replaced '$($r.old)' -> '$($r.new)': $found"
}
